$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the source diff (coin price/volume refresh).
$updates = @(
    @{Cell='D2'; Value='28.880.49'},
    @{Cell='E2'; Value='  -1.11%  '},
    @{Cell='D3'; Value='1.891.19'},
    @{Cell='E3'; Value='  -2.50%  '},
    @{Cell='D4'; Value='1.012'},
    @{Cell='E4'; Value='  +1.06%  '},
    @{Cell='D5'; Value='324.13'},
    @{Cell='E5'; Value='  -0.87%  '},
    @{Cell='D6'; Value='1.009'},
    @{Cell='E6'; Value='  +0.67%  '},
    @{Cell='D7'; Value='0.4582'},
    @{Cell='E7'; Value='  -1.16%  '},
    @{Cell='D8'; Value='0.3805'},
    @{Cell='E8'; Value='  -2.58%  '},
    @{Cell='B9'; Value='Dogecoin'},
    @{Cell='C9'; Value='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'},
    @{Cell='D9'; Value='0.07687'},
    @{Cell='E9'; Value='  -2.47%  '},
    @{Cell='B10'; Value='Polygon'},
    @{Cell='C10'; Value='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'},
    @{Cell='D10'; Value='0.9632'},
    @{Cell='E10'; Value='  -3.53%  '},
    @{Cell='B11'; Value='Solana'},
    @{Cell='C11'; Value='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'},
    @{Cell='D11'; Value='21.92'},
    @{Cell='E11'; Value='  -1.27%  '},
    @{Cell='B12'; Value='WrappedEther'},
    @{Cell='C12'; Value='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'},
    @{Cell='D12'; Value='1.897.27'},
    @{Cell='E12'; Value='  -1.63%  '},
    @{Cell='B13'; Value='Chainlink'},
    @{Cell='C13'; Value='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'},
    @{Cell='D13'; Value='6.964'},
    @{Cell='E13'; Value='  -2.01%  '},
    @{Cell='B14'; Value='Polkadot'},
    @{Cell='C14'; Value='https://coinranking.com/coin/25W7FG7om+polkadot-dot'},
    @{Cell='D14'; Value='5.674'},
    @{Cell='E14'; Value='  -2.63%  '},
    @{Cell='B15'; Value='TRON'},
    @{Cell='C15'; Value='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'},
    @{Cell='D15'; Value='0.07069'},
    @{Cell='E15'; Value='  -0.01%  '},
    @{Cell='B16'; Value='BinanceUSD'},
    @{Cell='C16'; Value='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'},
    @{Cell='D16'; Value='1.012'},
    @{Cell='E16'; Value='  +0.81%  '},
    @{Cell='B17'; Value='Litecoin'},
    @{Cell='C17'; Value='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'},
    @{Cell='D17'; Value='83.27'},
    @{Cell='E17'; Value='  -5.43%  '},
    @{Cell='B18'; Value='ShibaInu'},
    @{Cell='C18'; Value='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'},
    @{Cell='D18'; Value='0.000009407'},
    @{Cell='E18'; Value='  -5.71%  '},
    @{Cell='B19'; Value='Avalanche'},
    @{Cell='C19'; Value='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'},
    @{Cell='D19'; Value='16.79'},
    @{Cell='E19'; Value='  -2.28%  '},
    @{Cell='B20'; Value='Dai'},
    @{Cell='C20'; Value='https://coinranking.com/coin/MoTuySvg7+dai-dai'},
    @{Cell='D20'; Value='1.004'},
    @{Cell='E20'; Value='  +0.15%  '},
    @{Cell='B21'; Value='WrappedBTC'},
    @{Cell='C21'; Value='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'},
    @{Cell='D21'; Value='28.898.53'},
    @{Cell='E21'; Value='  -1.23%  '},
    @{Cell='B22'; Value='Uniswap'},
    @{Cell='C22'; Value='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'},
    @{Cell='D22'; Value='5.389'},
    @{Cell='E22'; Value='  -1.96%  '},
    @{Cell='B23'; Value='Cosmos'},
    @{Cell='C23'; Value='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'},
    @{Cell='D23'; Value='10.92'},
    @{Cell='E23'; Value='  -2.78%  '},
    @{Cell='B24'; Value='WrappedliquidstakedEther2.0'},
    @{Cell='C24'; Value='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'},
    @{Cell='D24'; Value='2.141.02'},
    @{Cell='E24'; Value='  -0.87%  '},
    @{Cell='B25'; Value='Toncoin'},
    @{Cell='C25'; Value='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'},
    @{Cell='D25'; Value='2.086'},
    @{Cell='E25'; Value='  -0.75%  '},
    @{Cell='B26'; Value='Monero'},
    @{Cell='C26'; Value='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'},
    @{Cell='D26'; Value='157.09'},
    @{Cell='E26'; Value='  +0.38%  '},
    @{Cell='B27'; Value='EthereumClassic'},
    @{Cell='C27'; Value='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'},
    @{Cell='D27'; Value='19.08'},
    @{Cell='E27'; Value='  -2.29%  '},
    @{Cell='B28'; Value='InternetComputer(DFINITY)'},
    @{Cell='C28'; Value='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'},
    @{Cell='D28'; Value='5.626'},
    @{Cell='E28'; Value='  -4.88%  '},
    @{Cell='B29'; Value='BitcoinCash'},
    @{Cell='C29'; Value='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'},
    @{Cell='D29'; Value='116.89'},
    @{Cell='E29'; Value='  -1.84%  '},
    @{Cell='B30'; Value='LidoDAOToken'},
    @{Cell='C30'; Value='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'},
    @{Cell='D30'; Value='1.810'},
    @{Cell='E30'; Value='  -3.95%  '},
    @{Cell='B31'; Value='Stellar'},
    @{Cell='C31'; Value='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'},
    @{Cell='D31'; Value='0.09290'},
    @{Cell='E31'; Value='  -0.69%  '},
    @{Cell='B32'; Value='ImmutableX'},
    @{Cell='C32'; Value='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'},
    @{Cell='D32'; Value='0.8513'},
    @{Cell='E32'; Value='  -4.85%  '},
    @{Cell='B33'; Value='Filecoin'},
    @{Cell='C33'; Value='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'},
    @{Cell='D33'; Value='5.059'},
    @{Cell='E33'; Value='  -3.41%  '},
    @{Cell='B34'; Value='ARBITRUM'},
    @{Cell='C34'; Value='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'},
    @{Cell='D34'; Value='1.246'},
    @{Cell='E34'; Value='  -5.97%  '},
    @{Cell='B35'; Value='HuobiToken'},
    @{Cell='C35'; Value='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'},
    @{Cell='D35'; Value='3.087'},
    @{Cell='E35'; Value='  -1.83%  '},
    @{Cell='B36'; Value='TrustWalletToken'},
    @{Cell='C36'; Value='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'},
    @{Cell='D36'; Value='1.160'},
    @{Cell='E36'; Value='  -1.21%  '},
    @{Cell='B37'; Value='Hedera'},
    @{Cell='C37'; Value='https://coinranking.com/coin/jad286TjB+hedera-hbar'},
    @{Cell='D37'; Value='0.05643'},
    @{Cell='E37'; Value='  -2.56%  '},
    @{Cell='B38'; Value='VeChain'},
    @{Cell='C38'; Value='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'},
    @{Cell='D38'; Value='0.02035'},
    @{Cell='E38'; Value='  -3.40%  '},
    @{Cell='B39'; Value='TheSandbox'},
    @{Cell='C39'; Value='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'},
    @{Cell='D39'; Value='0.5524'},
    @{Cell='E39'; Value='  -3.47%  '},
    @{Cell='D40'; Value='7.430'},
    @{Cell='E40'; Value='  -3.67%  '},
    @{Cell='B41'; Value='Algorand'},
    @{Cell='C41'; Value='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'},
    @{Cell='D41'; Value='0.1748'},
    @{Cell='E41'; Value='  -3.99%  '},
    @{Cell='B42'; Value='PEPE'},
    @{Cell='C42'; Value='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'},
    @{Cell='D42'; Value='0.000002839'},
    @{Cell='E42'; Value='  -10.23%  '},
    @{Cell='B43'; Value='Aptos'},
    @{Cell='C43'; Value='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'},
    @{Cell='D43'; Value='9.183'},
    @{Cell='E43'; Value='  -5.97%  '},
    @{Cell='B44'; Value='MXToken'},
    @{Cell='C44'; Value='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'},
    @{Cell='D44'; Value='2.714'},
    @{Cell='E44'; Value='  +4.39%  '},
    @{Cell='B45'; Value='Decentraland'},
    @{Cell='C45'; Value='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'},
    @{Cell='D45'; Value='0.5167'},
    @{Cell='E45'; Value='  -3.57%  '},
    @{Cell='B46'; Value='EnergySwap'},
    @{Cell='C46'; Value='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'},
    @{Cell='D46'; Value='11.25'},
    @{Cell='E46'; Value='  -5.59%  '},
    @{Cell='B47'; Value='Cronos'},
    @{Cell='C47'; Value='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'},
    @{Cell='D47'; Value='0.06784'},
    @{Cell='E47'; Value='  -2.40%  '},
    @{Cell='B48'; Value='RenderToken'},
    @{Cell='C48'; Value='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'},
    @{Cell='D48'; Value='2.058'},
    @{Cell='E48'; Value='  -6.95%  '},
    @{Cell='B49'; Value='NEARProtocol'},
    @{Cell='C49'; Value='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'},
    @{Cell='D49'; Value='1.774'},
    @{Cell='E49'; Value='  -4.38%  '},
    @{Cell='B50'; Value='Quant'},
    @{Cell='C50'; Value='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'},
    @{Cell='D50'; Value='110.09'},
    @{Cell='E50'; Value='  -2.80%  '},
    @{Cell='B51'; Value='WOONetwork'},
    @{Cell='C51'; Value='https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'},
    @{Cell='D51'; Value='0.2947'},
    @{Cell='E51'; Value='  -2.33%  '}
)

foreach ($u in $updates) {
    # Force text storage so numeric-looking strings (e.g. "1.012", "29.060.08")
    # are preserved verbatim instead of being parsed as numbers/dates.
    $ws.Range($u.Cell).NumberFormat = "@"
    $ws.Range($u.Cell).Value = $u.Value
    # Reset formatting back to the default style so no stray number format lingers.
    $ws.Range($u.Cell).Style = "Normal"
}
